# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to match the refreshed output data.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for the "展览" worksheet
$exhibitionUpdates = @{
    2  = 10608
    3  = 295
    6  = 802
    9  = 1143
    10 = 3317
    11 = 2454
    13 = 2248
    14 = 244
    18 = 608
    19 = 80
    20 = 271
    22 = 35
    24 = 55
    28 = 430
    29 = 626
    30 = 57
    31 = 62
    32 = 313
    34 = 1590
    35 = 631
    36 = 627
    37 = 1805
    38 = 176
    39 = 475
    41 = 502
    42 = 1118
    44 = 383
}

# Row -> new F-column value for the "全部类型" worksheet
$allTypesUpdates = @{
    2  = 10608
    3  = 295
    6  = 802
    8  = 1143
    9  = 3317
    10 = 2454
    11 = 2248
    12 = 244
    16 = 608
    17 = 80
    18 = 271
    20 = 35
    22 = 55
    26 = 430
    27 = 626
    28 = 57
    32 = 62
    33 = 313
    35 = 1590
    36 = 631
    38 = 627
    39 = 1805
    40 = 176
    44 = 475
    46 = 502
    47 = 1118
    48 = 383
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
